$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 "Save" - mirror the style of the existing header (G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data cells H2 / H3 = 0 (save/era data)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
